$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "Facility utilisation" header (was "Facility utilisatin")
$ws.Range("F1").Value = "Facility utilisation"

# Correct the utilisation percentage value for the "Fish GEF" facility row
$ws.Range("G4").Value = 124.75

# Header row is now one line shorter after the text fix; match the
# recalculated row height
$ws.Rows("1").RowHeight = 87.5

# Update the active selection to reflect where the user left off editing
$ws.Range("G14").Select()

Write-Host "Edit applied"
